# Employee_Data.xlsx update:
#  - drop the "is_admin" column and the bogus numeric "branch_num" column
#  - add a real "branch_num" column (employee number 1005151, integer format)
#  - collapse the duplicated "hire_date" column into a single column, leaving
#    the old duplicate as an empty (but still date-formatted) column
#  - fix the "Branch Manager" job title (it had a stray trailing quote)
#  - assorted cosmetic window/selection/column-width tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural column changes -------------------------------------------
# Column E ("is_admin") goes away entirely.
$ws.Columns.Item(5).Delete()
# The next column (the bogus numeric "branch_num" placeholder, now shifted
# into E) goes away too.
$ws.Columns.Item(5).Delete()
# Make room for the real "branch_num" values in front of the (now
# duplicated) hire_date columns.
$ws.Columns.Item(5).Insert()

# --- headers ---------------------------------------------------------------
$ws.Range("E1").Value = "branch_num"
$ws.Range("F1").Value = "hire_date"
# The extra duplicated "hire_date" header cell in G1 is no longer needed.
$ws.Range("G1").ClearContents()

# --- data --------------------------------------------------------------
# New branch_num values (integer display format, matches numFmtId 1 / "0").
$ws.Range("E2:E16").Value = 1005151
$ws.Range("E2:E16").NumberFormat = "0"

# The old duplicate hire_date column (now G) keeps its date formatting but
# the value itself is cleared out.
$ws.Range("G2:G16").ClearContents()

# Fix the "Branch Manager'" typo -> "Branch Manager" (keep quote-prefix
# formatting by re-entering it with a leading apostrophe, same as Excel's
# own quote-prefix input convention).
$ws.Range("D2").Value = "'Branch Manager"

# --- cosmetic tweaks -------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 17.5
$ws.Range("F24").Select()

$w = $wb.Windows.Item(1)
$w.Left = 28680
$w.Top = -120
